$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values for row 8 (C8:F8) from 2 to 5
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5

# Update the active cell/selection to G8
$ws.Range("G8").Select()
